$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as scraped by the GitHub Actions job.
# D-column "Price" values are forced to text (leading apostrophe) to avoid Excel
# auto-converting numeric-looking strings (e.g. "479.06") into numbers; the Style
# reset afterwards clears the resulting quote-prefix style so formatting matches
# the original (unstyled) cells.

$cell = $ws.Range("D2")
$cell.Value = "'56.021.96"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "
$cell = $ws.Range("D3")
$cell.Value = "'2.389.60"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.Value = "'479.06"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$cell = $ws.Range("D6")
$cell.Value = "'147.35"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.16%  "
$cell = $ws.Range("D8")
$cell.Value = "'0.499"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.44%  "
$cell = $ws.Range("D9")
$cell.Value = "'2.391.90"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -4.32%  "
$cell = $ws.Range("D10")
$cell.Value = "'0.0976"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$cell = $ws.Range("D11")
$cell.Value = "'5.43"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  +1.21%  "
$cell = $ws.Range("D14")
$cell.Value = "'2.809.36"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.63%  "
$cell = $ws.Range("D15")
$cell.Value = "'56.361.08"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "
$cell = $ws.Range("D16")
$cell.Value = "'20.34"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("E17").Value = "  -2.07%  "
$cell = $ws.Range("D18")
$cell.Value = "'2.392.34"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -4.47%  "
$cell = $ws.Range("D19")
$cell.Value = "'4.48"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "
$cell = $ws.Range("D20")
$cell.Value = "'315.24"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -4.15%  "
$cell = $ws.Range("D22")
$cell.Value = "'0.998"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$cell = $ws.Range("D23")
$cell.Value = "'5.68"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.52%  "
$cell = $ws.Range("D24")
$cell.Value = "'56.84"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("E27").Value = "  -4.15%  "
$cell = $ws.Range("D28")
$cell.Value = "'2.498.04"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -4.25%  "
$cell = $ws.Range("D29")
$cell.Value = "'7.26"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  +0.02%  "
$cell = $ws.Range("D32")
$cell.Value = "'148.70"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.29%  "
$cell = $ws.Range("D33")
$cell.Value = "'17.93"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  -0.26%  "
$cell = $ws.Range("D35")
$cell.Value = "'4.99"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D37")
$cell.Value = "'3.59"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D38")
$cell.Value = "'0.845"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$cell = $ws.Range("D39")
$cell.Value = "'33.41"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("E44").Value = "  -3.64%  "
$cell = $ws.Range("D45")
$cell.Value = "'0.0945"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +4.88%  "
$cell = $ws.Range("D46")
$cell.Value = "'10.23"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D47")
$cell.Value = "'254.39"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D48")
$cell.Value = "'4.63"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D49")
$cell.Value = "'0.0224"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$cell = $ws.Range("D50")
$cell.Value = "'17.05"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$cell = $ws.Range("D51")
$cell.Value = "'1.775.59"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -7.57%  "
